$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line to add the EM (10) entry.
$d.Content.Find.Execute("Curso (semestre ideal): EB (5)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EM (10), EB (6)", 2)

# 2. Remove the trailing "Requisitos" heading and its "LOT2059 ..." bullet
#    paragraph (the last two paragraphs of the document body).
$count = $d.Paragraphs.Count
$first = $d.Paragraphs($count - 1)
$last = $d.Paragraphs($count)
$r = $d.Range($first.Range.Start, $last.Range.End)
$r.Delete()
